$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '40.794.56'
$ws.Range("E2").Value = '  -2.87%  '

# Row 3
$ws.Range("D3").Value = '2.135.58'
$ws.Range("E3").Value = '  -3.60%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.41%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.594'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.35%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '68.05'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.80%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.554'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.54%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -11.38%  '

# Row 11
$ws.Range("E11").Value = '  -7.74%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '52.94'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.90%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0987'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.87%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.51%  '

# Row 15
$ws.Range("D15").Value = '2.450.22'
$ws.Range("E15").Value = '  -3.84%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.75%  '

# Row 17
$ws.Range("D17").Value = '2.128.39'
$ws.Range("E17").Value = '  -3.58%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.764'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -8.60%  '

# Row 19
$ws.Range("D19").Value = '40.593.11'
$ws.Range("E19").Value = '  -3.12%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0971'
$ws.Range("E20").Value = '  -9.67%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.62%  '

# Row 22
$ws.Range("E22").Value = '  -8.97%  '

# Row 23
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -11.27%  '

# Row 24
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '220.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.16%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.23%  '

# Row 26
$ws.Range("E26").Value = '  -12.14%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -13.45%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -13.61%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.53%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.05%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '166.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.99%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.33'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.69%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.21%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0735'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.64%  '

# Row 35
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.118'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.77%  '

# Row 36
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.93'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -12.43%  '

# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0966'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -11.80%  '

# Row 38
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.96'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.35%  '

# Row 39
$ws.Range("E39").Value = '  -10.04%  '

# Row 40
$ws.Range("E40").Value = '  -5.80%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -18.42%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.30%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '55.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -15.18%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.181'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.65%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.65%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0937'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.86%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '94.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.06%  '

# Row 48
$ws.Range("E48").Value = '  -5.95%  '

# Row 49
$ws.Range("E49").Value = '  -7.27%  '

# Row 50
$ws.Range("E50").Value = '  -3.80%  '

# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.09'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -13.84%  '
